$wb = $excel.ActiveWorkbook

# Duplicate the last existing sheet (Bus_Makulu_r) to use as the template for
# the new Trailer1Axle_f sheet, then place the copy at the end of the tab strip.
$template = $wb.Worksheets.Item("Bus_Makulu_r")
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Trailer1Axle_f"

# Update the hardpoint data for the new Trailer1Axle front droplink.
$new.Range("H3").Value = "Droplink_Trailer1Axle_f"

$new.Range("F5").Value = 0.05
$new.Range("G5").Value = 0.6
$new.Range("H5").Value = 0.19

$new.Range("F6").Formula = "=0.3-0.15"
$new.Range("G6").Value = 0.57999999999999996
$new.Range("H6").Value = 0.2

# Select/focus the new sheet as the active tab, matching the author's last
# saved view state.
$new.Activate()
$new.Range("H7").Select()
